$d = $word.ActiveDocument

function Replace-FirstMatch($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Replacements are applied in document order (top to bottom) so that a
# later "new" value that happens to equal an earlier "old" value does not
# accidentally get re-matched.
Replace-FirstMatch "2023-08-25 Friday" "2023-08-26 Saturday"

Replace-FirstMatch "42÷5=8, 2" "64÷9=7, 1"
Replace-FirstMatch "43÷7=6, 1" "12÷8=1, 4"
Replace-FirstMatch "46÷5=9, 1" "82÷4=20, 2"
Replace-FirstMatch "10÷6=1, 4" "57÷9=6, 3"
Replace-FirstMatch "11÷8=1, 3" "79÷7=11, 2"

Replace-FirstMatch "60÷5=12, 0" "93÷2=46, 1"
Replace-FirstMatch "31÷7=4, 3" "52÷3=17, 1"
Replace-FirstMatch "73÷6=12, 1" "38÷6=6, 2"
Replace-FirstMatch "26÷4=6, 2" "61÷3=20, 1"
Replace-FirstMatch "37÷3=12, 1" "55÷6=9, 1"

Replace-FirstMatch "29÷5=5, 4" "81÷8=10, 1"
Replace-FirstMatch "42÷6=7, 0" "70÷4=17, 2"
Replace-FirstMatch "29÷3=9, 2" "52÷8=6, 4"
Replace-FirstMatch "78÷6=13, 0" "66÷7=9, 3"
Replace-FirstMatch "21÷6=3, 3" "75÷2=37, 1"

Replace-FirstMatch "34÷9=3, 7" "41÷5=8, 1"
Replace-FirstMatch "75÷7=10, 5" "42÷5=8, 2"
Replace-FirstMatch "38÷7=5, 3" "93÷7=13, 2"
Replace-FirstMatch "25÷6=4, 1" "73÷7=10, 3"
Replace-FirstMatch "69÷6=11, 3" "17÷8=2, 1"

Replace-FirstMatch "56÷2=28, 0" "62÷4=15, 2"
Replace-FirstMatch "46÷3=15, 1" "85÷5=17, 0"
Replace-FirstMatch "36÷9=4, 0" "50÷4=12, 2"
Replace-FirstMatch "40÷4=10, 0" "34÷2=17, 0"
Replace-FirstMatch "96÷8=12, 0" "43÷3=14, 1"
